$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# B4: update path text
$ws.Range("B4").Value = "C:\Users\dpere\Documents\JTMT\forecast_by_version\V4\BASE_YEAR"

# B5: change from text "True" to boolean FALSE
$ws.Range("B5").Value = $false

# B6: clear the cell content (becomes empty inline string)
$ws.Range("B6").Value = ""
